$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.099.22'
$ws.Range("E2").Value = '  -0.67%  '
$ws.Range("D3").Value = '3.416.97'
$ws.Range("E3").Value = '  -1.44%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '407.56'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.06'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.56%  '
$ws.Range("E7").Value = '  -0.62%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.684'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.88%  '
$ws.Range("E10").Value = '  -5.47%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.79'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.81%  '
$ws.Range("E13").Value = '  -3.59%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.89'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.54%  '
$ws.Range("D15").Value = '3.431.98'
$ws.Range("E15").Value = '  -0.05%  '
$ws.Range("D16").Value = '62.132.63'
$ws.Range("E16").Value = '  -0.50%  '
$ws.Range("E17").Value = '  -3.08%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.03'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.97%  '
$ws.Range("E19").Value = '  -4.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.19'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.11%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '84.19'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '313.31'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.26%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.87'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.40%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.17'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.90%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.79'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +9.78%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '29.64'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.42%  '
$ws.Range("E27").Value = '  +0.06%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.79'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.18%  '
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.61'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.62%  '
$ws.Range("E30").Value = '  -2.33%  '
$ws.Range("E31").Value = '  -4.10%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '42.91'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.15%  '
$ws.Range("E33").Value = '  -0.26%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.39'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.99%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0483'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.29%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '51.89'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.50%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.36%  '
$ws.Range("E38").Value = '  -3.82%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.96'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.20%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.99'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.08%  '
$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.125'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.62%  '
$ws.Range("B42").Value = 'Monero'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '137.18'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.51%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.299'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.15%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.04'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.59%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.78'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.36%  '
$ws.Range("E46").Value = '  -2.38%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '21.48'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.75%  '
$ws.Range("D48").Value = '2.122.52'
$ws.Range("E48").Value = '  -4.55%  '
$ws.Range("E49").Value = '  -3.87%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.93'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.64'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +16.36%  '
